$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 253   # ALC!H55
$ws.Cells.Item(55, 9).Value = 248   # ALC!I55
$ws.Cells.Item(55, 10).Value = 260.5   # ALC!J55
$ws.Cells.Item(55, 11).Value = 248   # ALC!K55
$ws.Cells.Item(55, 12).Value = 260.5   # ALC!L55
$ws.Cells.Item(55, 13).Value = -34   # ALC!M55
$ws.Cells.Item(55, 14).Value = -688.5   # ALC!N55
$ws.Cells.Item(62, 8).Value = 2714.2856   # ALC!H62
$ws.Cells.Item(62, 9).Value = 2300.3333   # ALC!I62
$ws.Cells.Item(62, 11).Value = 2300.3333   # ALC!K62
$ws.Cells.Item(62, 13).Value = -1676.3333   # ALC!M62
$ws.Cells.Item(65, 8).Value = 2714.2856   # ALC!H65
$ws.Cells.Item(65, 9).Value = 2300.3333   # ALC!I65
$ws.Cells.Item(65, 11).Value = 11501.6665   # ALC!K65
$ws.Cells.Item(65, 13).Value = -8381.666499999999   # ALC!M65
$ws.Cells.Item(112, 8).Value = 2229.24   # ALC!H112
$ws.Cells.Item(112, 10).Value = 2349.1738   # ALC!J112
$ws.Cells.Item(112, 12).Value = 7047.5214   # ALC!L112
$ws.Cells.Item(112, 14).Value = -9263.5214   # ALC!N112
$ws.Cells.Item(123, 8).Value = 70699.09   # ALC!H123
$ws.Cells.Item(123, 10).Value = 70699.09   # ALC!J123
$ws.Cells.Item(123, 12).Value = 70699.09   # ALC!L123
$ws.Cells.Item(123, 14).Value = -80499.09   # ALC!N123
$ws.Cells.Item(135, 8).Value = 867.2   # ALC!H135
$ws.Cells.Item(135, 9).Value = 613.1724   # ALC!I135
$ws.Cells.Item(135, 10).Value = 2095   # ALC!J135
$ws.Cells.Item(135, 11).Value = 5518.551600000001   # ALC!K135
$ws.Cells.Item(135, 12).Value = 18855   # ALC!L135
$ws.Cells.Item(135, 13).Value = -2983.551600000001   # ALC!M135
$ws.Cells.Item(135, 14).Value = -23925   # ALC!N135
$ws.Cells.Item(137, 8).Value = 1705.4   # ALC!H137
$ws.Cells.Item(137, 9).Value = 1919.5555   # ALC!I137
$ws.Cells.Item(137, 10).Value = 1530.1818   # ALC!J137
$ws.Cells.Item(137, 11).Value = 5758.666499999999   # ALC!K137
$ws.Cells.Item(137, 12).Value = 4590.5454   # ALC!L137
$ws.Cells.Item(137, 13).Value = -3208.666499999999   # ALC!M137
$ws.Cells.Item(137, 14).Value = -9690.545399999999   # ALC!N137
$ws.Cells.Item(138, 8).Value = 4352784   # ALC!H138
$ws.Cells.Item(138, 9).Value = 9092613   # ALC!I138
$ws.Cells.Item(138, 10).Value = 7940.625   # ALC!J138
$ws.Cells.Item(138, 11).Value = 27277839   # ALC!K138
$ws.Cells.Item(138, 12).Value = 23821.875   # ALC!L138
$ws.Cells.Item(138, 13).Value = -27272699   # ALC!M138
$ws.Cells.Item(138, 14).Value = -34101.875   # ALC!N138
$ws.Cells.Item(141, 8).Value = 6956.2812   # ALC!H141
$ws.Cells.Item(141, 9).Value = 3954.8708   # ALC!I141
$ws.Cells.Item(141, 11).Value = 11864.6124   # ALC!K141
$ws.Cells.Item(141, 13).Value = -6684.6124   # ALC!M141

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 16450.654   # ARM!H32
$ws.Cells.Item(32, 9).Value = 16711.361   # ARM!I32
$ws.Cells.Item(32, 11).Value = 16711.361   # ARM!K32
$ws.Cells.Item(32, 13).Value = -16424.361   # ARM!M32
$ws.Cells.Item(45, 8).Value = 2227.75   # ARM!H45
$ws.Cells.Item(45, 9).Value = 2117.4285   # ARM!I45
$ws.Cells.Item(45, 11).Value = 2117.4285   # ARM!K45
$ws.Cells.Item(45, 13).Value = -1740.4285   # ARM!M45
$ws.Cells.Item(61, 8).Value = 1967.8889   # ARM!H61
$ws.Cells.Item(61, 9).Value = 1967.8889   # ARM!I61
$ws.Cells.Item(61, 11).Value = 1967.8889   # ARM!K61
$ws.Cells.Item(61, 13).Value = -1755.8889   # ARM!M61
$ws.Cells.Item(74, 8).Value = 1302.8235   # ARM!H74
$ws.Cells.Item(74, 9).Value = 1193.5   # ARM!I74
$ws.Cells.Item(74, 10).Value = 1400   # ARM!J74
$ws.Cells.Item(74, 11).Value = 1193.5   # ARM!K74
$ws.Cells.Item(74, 12).Value = 1400   # ARM!L74
$ws.Cells.Item(74, 13).Value = -319.5   # ARM!M74
$ws.Cells.Item(74, 14).Value = -3148   # ARM!N74
$ws.Cells.Item(77, 8).Value = 1302.8235   # ARM!H77
$ws.Cells.Item(77, 9).Value = 1193.5   # ARM!I77
$ws.Cells.Item(77, 10).Value = 1400   # ARM!J77
$ws.Cells.Item(77, 11).Value = 5967.5   # ARM!K77
$ws.Cells.Item(77, 12).Value = 7000   # ARM!L77
$ws.Cells.Item(77, 13).Value = -1599.5   # ARM!M77
$ws.Cells.Item(77, 14).Value = -15736   # ARM!N77
$ws.Cells.Item(122, 8).Value = 1815.909   # ARM!H122
$ws.Cells.Item(122, 9).Value = 1914   # ARM!I122
$ws.Cells.Item(122, 10).Value = 1482.4   # ARM!J122
$ws.Cells.Item(122, 11).Value = 5742   # ARM!K122
$ws.Cells.Item(122, 12).Value = 4447.200000000001   # ARM!L122
$ws.Cells.Item(122, 13).Value = -3292   # ARM!M122
$ws.Cells.Item(122, 14).Value = -9347.200000000001   # ARM!N122
$ws.Cells.Item(123, 8).Value = 43429   # ARM!H123
$ws.Cells.Item(123, 10).Value = 43429   # ARM!J123
$ws.Cells.Item(123, 12).Value = 43429   # ARM!L123
$ws.Cells.Item(123, 14).Value = -53229   # ARM!N123
$ws.Cells.Item(136, 8).Value = 1967.8889   # ARM!H136
$ws.Cells.Item(136, 9).Value = 1967.8889   # ARM!I136
$ws.Cells.Item(136, 11).Value = 5903.6667   # ARM!K136
$ws.Cells.Item(136, 13).Value = -3353.6667   # ARM!M136

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 24428.584   # BSM!H107
$ws.Cells.Item(107, 9).Value = 26322.092   # BSM!I107
$ws.Cells.Item(107, 10).Value = 3600   # BSM!J107
$ws.Cells.Item(107, 11).Value = 26322.092   # BSM!K107
$ws.Cells.Item(107, 12).Value = 3600   # BSM!L107
$ws.Cells.Item(107, 13).Value = -24402.092   # BSM!M107
$ws.Cells.Item(107, 14).Value = -7440   # BSM!N107

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 5000   # CRP!H50
$ws.Cells.Item(50, 10).Value = 0   # CRP!J50
$ws.Cells.Item(50, 12).Value = 0   # CRP!L50
$ws.Cells.Item(50, 14).Value = $null   # CRP!N50 remove
$ws.Cells.Item(99, 8).Value = 2085.087   # CRP!H99
$ws.Cells.Item(99, 9).Value = 2350.1765   # CRP!I99
$ws.Cells.Item(99, 11).Value = 2350.1765   # CRP!K99
$ws.Cells.Item(99, 13).Value = -852.1765   # CRP!M99
$ws.Cells.Item(122, 8).Value = 1106.0555   # CRP!H122
$ws.Cells.Item(122, 9).Value = 1199.6923   # CRP!I122
$ws.Cells.Item(122, 10).Value = 862.6   # CRP!J122
$ws.Cells.Item(122, 11).Value = 3599.0769   # CRP!K122
$ws.Cells.Item(122, 12).Value = 2587.8   # CRP!L122
$ws.Cells.Item(122, 13).Value = -1149.0769   # CRP!M122
$ws.Cells.Item(122, 14).Value = -7487.8   # CRP!N122
$ws.Cells.Item(126, 8).Value = 2085.087   # CRP!H126
$ws.Cells.Item(126, 9).Value = 2350.1765   # CRP!I126
$ws.Cells.Item(126, 11).Value = 7050.529500000001   # CRP!K126
$ws.Cells.Item(126, 13).Value = -4580.529500000001   # CRP!M126

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 94.333336   # CUL!H4
$ws.Cells.Item(4, 9).Value = 94.333336   # CUL!I4
$ws.Cells.Item(4, 10).Value = 0   # CUL!J4
$ws.Cells.Item(4, 11).Value = 283.000008   # CUL!K4
$ws.Cells.Item(4, 12).Value = 0   # CUL!L4
$ws.Cells.Item(4, 13).Value = -171.000008   # CUL!M4
$ws.Cells.Item(4, 14).Value = $null   # CUL!N4 remove
$ws.Cells.Item(61, 8).Value = 217.66667   # CUL!H61
$ws.Cells.Item(61, 9).Value = 76.5   # CUL!I61
$ws.Cells.Item(61, 10).Value = 500   # CUL!J61
$ws.Cells.Item(61, 11).Value = 229.5   # CUL!K61
$ws.Cells.Item(61, 12).Value = 1500   # CUL!L61
$ws.Cells.Item(61, 13).Value = -14.5   # CUL!M61
$ws.Cells.Item(61, 14).Value = -1930   # CUL!N61
$ws.Cells.Item(68, 8).Value = 1436.8667   # CUL!H68
$ws.Cells.Item(68, 10).Value = 1350.2307   # CUL!J68
$ws.Cells.Item(68, 12).Value = 4050.6921   # CUL!L68
$ws.Cells.Item(68, 14).Value = -5672.6921   # CUL!N68
$ws.Cells.Item(71, 8).Value = 1436.8667   # CUL!H71
$ws.Cells.Item(71, 10).Value = 1350.2307   # CUL!J71
$ws.Cells.Item(71, 12).Value = 12152.0763   # CUL!L71
$ws.Cells.Item(71, 14).Value = -20264.0763   # CUL!N71
$ws.Cells.Item(132, 8).Value = 2291.6   # CUL!H132
$ws.Cells.Item(132, 9).Value = 1250.5   # CUL!I132
$ws.Cells.Item(132, 10).Value = 2551.875   # CUL!J132
$ws.Cells.Item(132, 11).Value = 11254.5   # CUL!K132
$ws.Cells.Item(132, 12).Value = 22966.875   # CUL!L132
$ws.Cells.Item(132, 13).Value = -8724.5   # CUL!M132
$ws.Cells.Item(132, 14).Value = -28026.875   # CUL!N132

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 140000   # GSM!H15
$ws.Cells.Item(15, 10).Value = 140000   # GSM!J15
$ws.Cells.Item(15, 12).Value = 140000   # GSM!L15
$ws.Cells.Item(15, 14).Value = -140576   # GSM!N15
$ws.Cells.Item(43, 8).Value = 3057   # GSM!H43
$ws.Cells.Item(43, 9).Value = 1585.5   # GSM!I43
$ws.Cells.Item(43, 11).Value = 1585.5   # GSM!K43
$ws.Cells.Item(43, 13).Value = -1434.5   # GSM!M43
$ws.Cells.Item(53, 8).Value = 12494.75   # GSM!H53
$ws.Cells.Item(53, 9).Value = 4999.5   # GSM!I53
$ws.Cells.Item(53, 10).Value = 19990   # GSM!J53
$ws.Cells.Item(53, 11).Value = 4999.5   # GSM!K53
$ws.Cells.Item(53, 12).Value = 19990   # GSM!L53
$ws.Cells.Item(53, 13).Value = -4368.5   # GSM!M53
$ws.Cells.Item(53, 14).Value = -21252   # GSM!N53
$ws.Cells.Item(63, 8).Value = 35000   # GSM!H63
$ws.Cells.Item(63, 10).Value = 35000   # GSM!J63
$ws.Cells.Item(63, 12).Value = 35000   # GSM!L63
$ws.Cells.Item(63, 14).Value = -36372   # GSM!N63
$ws.Cells.Item(66, 8).Value = 35000   # GSM!H66
$ws.Cells.Item(66, 10).Value = 35000   # GSM!J66
$ws.Cells.Item(66, 12).Value = 105000   # GSM!L66
$ws.Cells.Item(66, 14).Value = -111864   # GSM!N66
$ws.Cells.Item(80, 8).Value = 3177.6155   # GSM!H80
$ws.Cells.Item(80, 9).Value = 3600.5   # GSM!I80
$ws.Cells.Item(80, 11).Value = 3600.5   # GSM!K80
$ws.Cells.Item(80, 13).Value = -2602.5   # GSM!M80
$ws.Cells.Item(81, 8).Value = 140000   # GSM!H81
$ws.Cells.Item(81, 10).Value = 140000   # GSM!J81
$ws.Cells.Item(81, 12).Value = 140000   # GSM!L81
$ws.Cells.Item(81, 14).Value = -141996   # GSM!N81
$ws.Cells.Item(83, 8).Value = 3177.6155   # GSM!H83
$ws.Cells.Item(83, 9).Value = 3600.5   # GSM!I83
$ws.Cells.Item(83, 11).Value = 18002.5   # GSM!K83
$ws.Cells.Item(83, 13).Value = -13010.5   # GSM!M83
$ws.Cells.Item(84, 8).Value = 140000   # GSM!H84
$ws.Cells.Item(84, 10).Value = 140000   # GSM!J84
$ws.Cells.Item(84, 12).Value = 420000   # GSM!L84
$ws.Cells.Item(84, 14).Value = -429984   # GSM!N84
$ws.Cells.Item(109, 8).Value = 12280.1   # GSM!H109
$ws.Cells.Item(109, 10).Value = 12280.1   # GSM!J109
$ws.Cells.Item(109, 12).Value = 12280.1   # GSM!L109
$ws.Cells.Item(109, 14).Value = -14360.1   # GSM!N109
$ws.Cells.Item(122, 8).Value = 3718   # GSM!H122
$ws.Cells.Item(122, 9).Value = 4515.846   # GSM!I122
$ws.Cells.Item(122, 10).Value = 2565.5557   # GSM!J122
$ws.Cells.Item(122, 11).Value = 13547.538   # GSM!K122
$ws.Cells.Item(122, 12).Value = 7696.6671   # GSM!L122
$ws.Cells.Item(122, 13).Value = -11097.538   # GSM!M122
$ws.Cells.Item(122, 14).Value = -12596.6671   # GSM!N122
$ws.Cells.Item(132, 8).Value = 2195.0908   # GSM!H132
$ws.Cells.Item(132, 9).Value = 1538.5   # GSM!I132
$ws.Cells.Item(132, 10).Value = 3205.2307   # GSM!J132
$ws.Cells.Item(132, 11).Value = 4615.5   # GSM!K132
$ws.Cells.Item(132, 12).Value = 9615.6921   # GSM!L132
$ws.Cells.Item(132, 13).Value = -2085.5   # GSM!M132
$ws.Cells.Item(132, 14).Value = -14675.6921   # GSM!N132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 11996.667   # LTW!H40
$ws.Cells.Item(40, 9).Value = 16495   # LTW!I40
$ws.Cells.Item(40, 10).Value = 3000   # LTW!J40
$ws.Cells.Item(40, 11).Value = 16495   # LTW!K40
$ws.Cells.Item(40, 12).Value = 3000   # LTW!L40
$ws.Cells.Item(40, 13).Value = -16359   # LTW!M40
$ws.Cells.Item(40, 14).Value = -3272   # LTW!N40

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 850   # WVR!H96
$ws.Cells.Item(96, 10).Value = 800   # WVR!J96
$ws.Cells.Item(96, 12).Value = 800   # WVR!L96
$ws.Cells.Item(96, 14).Value = -3546   # WVR!N96
$ws.Cells.Item(113, 8).Value = 734.375   # WVR!H113
$ws.Cells.Item(113, 9).Value = 453.4737   # WVR!I113
$ws.Cells.Item(113, 10).Value = 1801.8   # WVR!J113
$ws.Cells.Item(113, 11).Value = 1360.4211   # WVR!K113
$ws.Cells.Item(113, 12).Value = 5405.4   # WVR!L113
$ws.Cells.Item(113, 13).Value = 809.5789   # WVR!M113
$ws.Cells.Item(113, 14).Value = -9745.4   # WVR!N113
$ws.Cells.Item(122, 8).Value = 22729656   # WVR!H122
$ws.Cells.Item(122, 9).Value = 35715600   # WVR!I122
$ws.Cells.Item(122, 10).Value = 4252.5   # WVR!J122
$ws.Cells.Item(122, 11).Value = 107146800   # WVR!K122
$ws.Cells.Item(122, 12).Value = 12757.5   # WVR!L122
$ws.Cells.Item(122, 13).Value = -107144350   # WVR!M122
$ws.Cells.Item(122, 14).Value = -17657.5   # WVR!N122
